# Apply the data update described in the commit: extend the Bmp8a-Tgfbr1
# LR-pair analysis to include "sCs" as an additional Sending cluster (in
# addition to already being a Target cluster), producing a full 3x3 grid
# of Sending x Target clusters (ECs, FAPs, sCs) instead of the previous 2x3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bmp8a"
$ws.Cells.Item(2,3).Value = "Tgfbr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.3786073333333334
$ws.Cells.Item(2,8).Value = 1.135822
$ws.Cells.Item(2,9).Value = 0.3713290366620658
$ws.Cells.Item(2,10).Value = 0.3713290366620658
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 63.76294533333333
$ws.Cells.Item(2,14).Value = 191.288836
$ws.Cells.Item(2,15).Value = 0.6446527016991613
$ws.Cells.Item(2,16).Value = 0.6446527016991614
$ws.Cells.Item(2,17).Value = 24.14111869813244
$ws.Cells.Item(2,18).Value = 217.270068283192
$ws.Cells.Item(2,19).Value = 0.2393782667035476
$ws.Cells.Item(2,20).Value = 0.2393782667035477

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bmp8a"
$ws.Cells.Item(3,3).Value = "Tgfbr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.3786073333333334
$ws.Cells.Item(3,8).Value = 1.135822
$ws.Cells.Item(3,9).Value = 0.3713290366620658
$ws.Cells.Item(3,10).Value = 0.3713290366620658
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 22.32219333333333
$ws.Cells.Item(3,14).Value = 66.96658
$ws.Cells.Item(3,15).Value = 0.2256806388876402
$ws.Cells.Item(3,16).Value = 0.2256806388876402
$ws.Cells.Item(3,17).Value = 8.451346092084444
$ws.Cells.Item(3,18).Value = 76.06211482876
$ws.Cells.Item(3,19).Value = 0.08380177423142697
$ws.Cells.Item(3,20).Value = 0.08380177423142697

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bmp8a"
$ws.Cells.Item(4,3).Value = "Tgfbr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.3786073333333334
$ws.Cells.Item(4,8).Value = 1.135822
$ws.Cells.Item(4,9).Value = 0.3713290366620658
$ws.Cells.Item(4,10).Value = 0.3713290366620658
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 12.825399
$ws.Cells.Item(4,14).Value = 38.476197
$ws.Cells.Item(4,15).Value = 0.1296666594131984
$ws.Cells.Item(4,16).Value = 0.1296666594131984
$ws.Cells.Item(4,17).Value = 4.855790114326
$ws.Cells.Item(4,18).Value = 43.70211102893401
$ws.Cells.Item(4,19).Value = 0.04814899572709115
$ws.Cells.Item(4,20).Value = 0.04814899572709115

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Bmp8a"
$ws.Cells.Item(5,3).Value = "Tgfbr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.4265683333333333
$ws.Cells.Item(5,8).Value = 1.279705
$ws.Cells.Item(5,9).Value = 0.4183680408212104
$ws.Cells.Item(5,10).Value = 0.4183680408212104
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 63.76294533333333
$ws.Cells.Item(5,14).Value = 191.288836
$ws.Cells.Item(5,15).Value = 0.6446527016991613
$ws.Cells.Item(5,16).Value = 0.6446527016991614
$ws.Cells.Item(5,17).Value = 27.19925331926444
$ws.Cells.Item(5,18).Value = 244.79327987338
$ws.Cells.Item(5,19).Value = 0.2697020878199783
$ws.Cells.Item(5,20).Value = 0.2697020878199783

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Bmp8a"
$ws.Cells.Item(6,3).Value = "Tgfbr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.4265683333333333
$ws.Cells.Item(6,8).Value = 1.279705
$ws.Cells.Item(6,9).Value = 0.4183680408212104
$ws.Cells.Item(6,10).Value = 0.4183680408212104
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 22.32219333333333
$ws.Cells.Item(6,14).Value = 66.96658
$ws.Cells.Item(6,15).Value = 0.2256806388876402
$ws.Cells.Item(6,16).Value = 0.2256806388876402
$ws.Cells.Item(6,17).Value = 9.521940806544443
$ws.Cells.Item(6,18).Value = 85.69746725889998
$ws.Cells.Item(6,19).Value = 0.09441756674270108
$ws.Cells.Item(6,20).Value = 0.09441756674270108

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Bmp8a"
$ws.Cells.Item(7,3).Value = "Tgfbr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.4265683333333333
$ws.Cells.Item(7,8).Value = 1.279705
$ws.Cells.Item(7,9).Value = 0.4183680408212104
$ws.Cells.Item(7,10).Value = 0.4183680408212104
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 12.825399
$ws.Cells.Item(7,14).Value = 38.476197
$ws.Cells.Item(7,15).Value = 0.1296666594131984
$ws.Cells.Item(7,16).Value = 0.1296666594131984
$ws.Cells.Item(7,17).Value = 5.470909075764999
$ws.Cells.Item(7,18).Value = 49.238181681885
$ws.Cells.Item(7,19).Value = 0.05424838625853098
$ws.Cells.Item(7,20).Value = 0.05424838625853098

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Bmp8a"
$ws.Cells.Item(8,3).Value = "Tgfbr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.214425
$ws.Cells.Item(8,8).Value = 0.6432749999999999
$ws.Cells.Item(8,9).Value = 0.2103029225167239
$ws.Cells.Item(8,10).Value = 0.2103029225167238
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 63.76294533333333
$ws.Cells.Item(8,14).Value = 191.288836
$ws.Cells.Item(8,15).Value = 0.6446527016991613
$ws.Cells.Item(8,16).Value = 0.6446527016991614
$ws.Cells.Item(8,17).Value = 13.6723695531
$ws.Cells.Item(8,18).Value = 123.0513259779
$ws.Cells.Item(8,19).Value = 0.1355723471756354
$ws.Cells.Item(8,20).Value = 0.1355723471756354

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Bmp8a"
$ws.Cells.Item(9,3).Value = "Tgfbr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.214425
$ws.Cells.Item(9,8).Value = 0.6432749999999999
$ws.Cells.Item(9,9).Value = 0.2103029225167239
$ws.Cells.Item(9,10).Value = 0.2103029225167238
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 22.32219333333333
$ws.Cells.Item(9,14).Value = 66.96658
$ws.Cells.Item(9,15).Value = 0.2256806388876402
$ws.Cells.Item(9,16).Value = 0.2256806388876402
$ws.Cells.Item(9,17).Value = 4.786436305499999
$ws.Cells.Item(9,18).Value = 43.07792674949999
$ws.Cells.Item(9,19).Value = 0.04746129791351213
$ws.Cells.Item(9,20).Value = 0.04746129791351212

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Bmp8a"
$ws.Cells.Item(10,3).Value = "Tgfbr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.214425
$ws.Cells.Item(10,8).Value = 0.6432749999999999
$ws.Cells.Item(10,9).Value = 0.2103029225167239
$ws.Cells.Item(10,10).Value = 0.2103029225167238
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 12.825399
$ws.Cells.Item(10,14).Value = 38.476197
$ws.Cells.Item(10,15).Value = 0.1296666594131984
$ws.Cells.Item(10,16).Value = 0.1296666594131984
$ws.Cells.Item(10,17).Value = 2.750086180574999
$ws.Cells.Item(10,18).Value = 24.750775625175
$ws.Cells.Item(10,19).Value = 0.02726927742757629
$ws.Cells.Item(10,20).Value = 0.02726927742757628
